$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MongoDB")

# Rows 2-6 data (CmRDT-O style replication, now MongoDB numbers)
$ws.Range("A2").Value = 138.1543
$ws.Range("B2").Value = 100
$ws.Range("D2").Value = 141.57635
$ws.Range("E2").Value = 97.25536
$ws.Range("G2").Value = 149.01117
$ws.Range("H2").Value = 97.3747
$ws.Range("J2").Value = 216.31117
$ws.Range("K2").Value = 94.03342
$ws.Range("M2").Value = 544.25745
$ws.Range("N2").Value = 94.03342
$ws.Range("A3").Value = 136.91359
$ws.Range("B3").Value = 100
$ws.Range("D3").Value = 141.46983
$ws.Range("E3").Value = 97.0167
$ws.Range("G3").Value = 157.68929
$ws.Range("H3").Value = 98.09069
$ws.Range("J3").Value = 219.11932
$ws.Range("K3").Value = 97.7327
$ws.Range("M3").Value = 643.39124
$ws.Range("N3").Value = 95.584724
$ws.Range("A4").Value = 137.48409
$ws.Range("B4").Value = 100
$ws.Range("D4").Value = 140.36197
$ws.Range("E4").Value = 97.13604
$ws.Range("G4").Value = 156.97876
$ws.Range("H4").Value = 98.09069
$ws.Range("J4").Value = 287.93222
$ws.Range("K4").Value = 97.25536
$ws.Range("M4").Value = 634.1802
$ws.Range("N4").Value = 87.828156
$ws.Range("A5").Value = 136.80052
$ws.Range("B5").Value = 100
$ws.Range("D5").Value = 140.78091
$ws.Range("E5").Value = 96.65871
$ws.Range("G5").Value = 149.83243
$ws.Range("H5").Value = 98.09069
$ws.Range("J5").Value = 236.29297
$ws.Range("K5").Value = 92.24344
$ws.Range("M5").Value = 592.6044
$ws.Range("N5").Value = 95.1074
$ws.Range("A6").Value = 137.57307
$ws.Range("B6").Value = 100
$ws.Range("D6").Value = 139.66197
$ws.Range("E6").Value = 96.77803
$ws.Range("G6").Value = 148.81302
$ws.Range("H6").Value = 98.09069
$ws.Range("J6").Value = 257.21576
$ws.Range("K6").Value = 96.18139
$ws.Range("M6").Value = 692.1512
$ws.Range("N6").Value = 93.31744

# Rows 10-14 data
$ws.Range("A10").Value = 1379.0962
$ws.Range("B10").Value = 49.284008
$ws.Range("D10").Value = 1470.3536
$ws.Range("E10").Value = 64.916466
$ws.Range("G10").Value = 1635.4519
$ws.Range("H10").Value = 48.448692
$ws.Range("J10").Value = 1814.5366
$ws.Range("K10").Value = 49.164677
$ws.Range("M10").Value = 1919.6003
$ws.Range("N10").Value = 47.613365
$ws.Range("A11").Value = 1254.6715
$ws.Range("B11").Value = 48.56802
$ws.Range("D11").Value = 1377.5021
$ws.Range("E11").Value = 64.55847
$ws.Range("G11").Value = 1756.8237
$ws.Range("H11").Value = 51.193317
$ws.Range("J11").Value = 1906.9856
$ws.Range("K11").Value = 42.124104
$ws.Range("M11").Value = 1878.6364
$ws.Range("N11").Value = 47.613365
$ws.Range("A12").Value = 1249.1376
$ws.Range("B12").Value = 52.38663
$ws.Range("D12").Value = 1520.3763
$ws.Range("E12").Value = 41.408115
$ws.Range("G12").Value = 1764.0454
$ws.Range("H12").Value = 46.65872
$ws.Range("J12").Value = 1875.1229
$ws.Range("K12").Value = 48.44869
$ws.Range("M12").Value = 1891.1652
$ws.Range("N12").Value = 47.136044
$ws.Range("A13").Value = 1332.8982
$ws.Range("B13").Value = 47.732697
$ws.Range("D13").Value = 1505.4059
$ws.Range("E13").Value = 56.20525
$ws.Range("G13").Value = 1680.488
$ws.Range("H13").Value = 49.761333
$ws.Range("J13").Value = 1872.7358
$ws.Range("K13").Value = 49.403347
$ws.Range("M13").Value = 1967.0757
$ws.Range("N13").Value = 47.732697
$ws.Range("A14").Value = 1375.927
$ws.Range("B14").Value = 51.31265
$ws.Range("D14").Value = 1449.2662
$ws.Range("E14").Value = 65.632454
$ws.Range("G14").Value = 1723.487
$ws.Range("H14").Value = 46.65872
$ws.Range("J14").Value = 2006.7717
$ws.Range("K14").Value = 53.22196
$ws.Range("M14").Value = 1867.5139
$ws.Range("N14").Value = 48.806683

# Rows 18-22 data
$ws.Range("A18").Value = 1962.2244
$ws.Range("B18").Value = 52.625298
$ws.Range("D18").Value = 1903.2737
$ws.Range("E18").Value = 52.505962
$ws.Range("G18").Value = 1957.1813
$ws.Range("H18").Value = 52.2673
$ws.Range("A19").Value = 1992.4329
$ws.Range("B19").Value = 48.926014
$ws.Range("D19").Value = 1933.1799
$ws.Range("E19").Value = 47.971355
$ws.Range("G19").Value = 1983.0354
$ws.Range("H19").Value = 49.522675
$ws.Range("A20").Value = 1915.9952
$ws.Range("B20").Value = 49.642006
$ws.Range("D20").Value = 1892.4719
$ws.Range("E20").Value = 49.045345
$ws.Range("G20").Value = 1907.6597
$ws.Range("H20").Value = 50.835323
$ws.Range("A21").Value = 1896.6304
$ws.Range("B21").Value = 47.613365
$ws.Range("D21").Value = 1889.6527
$ws.Range("E21").Value = 46.897373
$ws.Range("G21").Value = 1922.0763
$ws.Range("H21").Value = 45.226727
$ws.Range("A22").Value = 1928.9751
$ws.Range("B22").Value = 54.415276
$ws.Range("D22").Value = 1955.5447
$ws.Range("E22").Value = 48.68735
$ws.Range("G22").Value = 1928.6742
$ws.Range("H22").Value = 55.608593

# Activate MongoDB sheet and set its selection, matching the new tabSelected state
$ws.Activate()
$ws.Range("A22").Select()
